# The "2567" sheet contained a duplicate/leaked row (row 6, "@Chompunoot ")
# that already exists derived via a CONCAT formula in row 5. Remove the
# redundant row entirely so every row below shifts up by one, matching the
# cleaned-up workbook (commit: "Removed secrets and cleaned history").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2567")
$ws.Activate()

# Select & delete the whole row 6 (shifts rows 7:52 up to 6:51).
$ws.Rows(6).Select() | Out-Null
$ws.Rows(6).Delete() | Out-Null
